$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
try {
  $d = $nm.Design
  Write-Output $d
  Write-Output $d.Index
  Write-Output $d.Name
} catch {
  Write-Output "ERR: $_"
}
